$p = $ppt.ActivePresentation
$p.Slides.Item(4).Delete()
